$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Lamp ("l" / "Lâmpada"), quantity 0, boolean flag FALSE
$ws.Cells.Item(2, 1).Value = "l"
$ws.Cells.Item(2, 2).Value = "Lâmpada"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = $false

# Row 3: Air conditioner ("a" / "A/C"), quantity 23, boolean flag FALSE
$ws.Cells.Item(3, 1).Value = "a"
$ws.Cells.Item(3, 2).Value = "A/C"
$ws.Cells.Item(3, 3).Value = 23
$ws.Cells.Item(3, 4).Value = $false
